{"js": "// Replace the date heading and the multiplication problems in the table\n// with the new values from the target revision. Each old string is\n// unique in the document, so a simple search-and-replace per pair is\n// safe and unambiguous.\nconst replacements = [\n  [\"2025-12-14 Sunday\", \"2025-12-15 Monday\"],\n  [\"65\u00d715=\", \"15\u00d733=\"],\n  [\"36\u00d712=\", \"57\u00d714=\"],\n  [\"37\u00d717=\", \"99\u00d714=\"],\n  [\"15\u00d793=\", \"14\u00d774=\"],\n  [\"51\u00d781=\", \"12\u00d762=\"],\n  [\"16\u00d740=\", \"18\u00d781=\"],\n  [\"19\u00d733=\", \"90\u00d772=\"],\n  [\"56\u00d765=\", \"46\u00d781=\"],\n  [\"35\u00d798=\", \"11\u00d759=\"],\n  [\"71\u00d744=\", \"67\u00d726=\"],\n  [\"65\u00d763=\", \"66\u00d777=\"],\n  [\"93\u00d719=\", \"19\u00d742=\"],\n  [\"90\u00d779=\", \"25\u00d776=\"],\n  [\"95\u00d721=\", \"85\u00d712=\"],\n  [\"94\u00d743=\", \"33\u00d776=\"],\n  [\"56\u00d739=\", \"58\u00d761=\"],\n  [\"59\u00d766=\", \"68\u00d745=\"],\n  [\"60\u00d754=\", \"92\u00d723=\"],\n  [\"34\u00d799=\", \"82\u00d717=\"],\n  [\"61\u00d733=\", \"37\u00d787=\"],\n  [\"34\u00d757=\", \"80\u00d722=\"],\n  [\"15\u00d714=\", \"50\u00d791=\"],\n  [\"13\u00d726=\", \"27\u00d799=\"],\n  [\"46\u00d798=\", \"30\u00d747=\"],\n  [\"72\u00d734=\", \"39\u00d734=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and the multiplication problems in the table\n# with the new values from the target revision. Each old string is\n# unique in the document, so a Find/Replace pass per pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-12-14 Sunday\", \"2025-12-15 Monday\"),\n  @(\"65\u00d715=\", \"15\u00d733=\"),\n  @(\"36\u00d712=\", \"57\u00d714=\"),\n  @(\"37\u00d717=\", \"99\u00d714=\"),\n  @(\"15\u00d793=\", \"14\u00d774=\"),\n  @(\"51\u00d781=\", \"12\u00d762=\"),\n  @(\"16\u00d740=\", \"18\u00d781=\"),\n  @(\"19\u00d733=\", \"90\u00d772=\"),\n  @(\"56\u00d765=\", \"46\u00d781=\"),\n  @(\"35\u00d798=\", \"11\u00d759=\"),\n  @(\"71\u00d744=\", \"67\u00d726=\"),\n  @(\"65\u00d763=\", \"66\u00d777=\"),\n  @(\"93\u00d719=\", \"19\u00d742=\"),\n  @(\"90\u00d779=\", \"25\u00d776=\"),\n  @(\"95\u00d721=\", \"85\u00d712=\"),\n  @(\"94\u00d743=\", \"33\u00d776=\"),\n  @(\"56\u00d739=\", \"58\u00d761=\"),\n  @(\"59\u00d766=\", \"68\u00d745=\"),\n  @(\"60\u00d754=\", \"92\u00d723=\"),\n  @(\"34\u00d799=\", \"82\u00d717=\"),\n  @(\"61\u00d733=\", \"37\u00d787=\"),\n  @(\"34\u00d757=\", \"80\u00d722=\"),\n  @(\"15\u00d714=\", \"50\u00d791=\"),\n  @(\"13\u00d726=\", \"27\u00d799=\"),\n  @(\"46\u00d798=\", \"30\u00d747=\"),\n  @(\"72\u00d734=\", \"39\u00d734=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
